$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = [double]"23.80000000000028"
$ws.Range("H2").Value = [double]"0.0003220541921479247"
$ws.Range("I2").Value = [double]"0.0003220541921479247"
$ws.Range("L2").Value = [double]"40.63506410647648"
$ws.Range("M2").Value = "[19.06661478302071, 62.203513429932244]"
$ws.Range("N2").Value = [double]"0.0004391639257395052"
$ws.Range("O2").Value = [double]"0.0004391639257395052"
$ws.Range("P2").Value = [double]"1.62897396852804"
$ws.Range("Q2").Value = "[0.9622896416401163, 2.295658295415964]"
$ws.Range("R2").Value = [double]"1.195143531074194e-05"
$ws.Range("S2").Value = [double]"1.195143531074194e-05"
$ws.Range("T2").Value = [double]"51.06925899217629"
$ws.Range("U2").Value = "[37.983612786938295, 64.15490519741428]"
$ws.Range("V2").Value = [double]"5.507010403249524e-10"
$ws.Range("W2").Value = [double]"5.507010403249524e-10"
$ws.Range("X2").Value = [double]"17.62962962962984"
$ws.Range("Y2").Value = [double]"15.10430430430448"
$ws.Range("Z2").Value = [double]"20.15495495495519"
$ws.Range("F3").Value = [double]"23.80000000000028"
$ws.Range("H3").Value = [double]"6.394189716707466e-06"
$ws.Range("I3").Value = [double]"6.394189716707466e-06"
$ws.Range("L3").Value = [double]"57.25997674926753"
$ws.Range("M3").Value = "[33.093457788081906, 81.42649571045315]"
$ws.Range("N3").Value = [double]"1.956131823721563e-05"
$ws.Range("O3").Value = [double]"1.956131823721563e-05"
$ws.Range("P3").Value = [double]"1.855395060678656"
$ws.Range("Q3").Value = "[1.3648160276856558, 2.345974093671657]"
$ws.Range("R3").Value = [double]"1.249431447192251e-09"
$ws.Range("S3").Value = [double]"1.249431447192251e-09"
$ws.Range("T3").Value = [double]"61.95128043618536"
$ws.Range("U3").Value = "[47.49736582059712, 76.4051950517736]"
$ws.Range("V3").Value = [double]"4.213718263201827e-11"
$ws.Range("W3").Value = [double]"4.213718263201827e-11"
$ws.Range("X3").Value = [double]"16.77197197197217"
$ws.Range("Y3").Value = [double]"14.91371371371389"
$ws.Range("Z3").Value = [double]"18.63023023023045"
$ws.Range("F4").Value = [double]"23.80000000000028"
$ws.Range("H4").Value = [double]"0.0002463665194947229"
$ws.Range("I4").Value = [double]"0.0002463665194947229"
$ws.Range("L4").Value = [double]"35.8293676700638"
$ws.Range("M4").Value = "[15.628524745311829, 56.03021059481576]"
$ws.Range("N4").Value = [double]"0.0008571597498645556"
$ws.Range("O4").Value = [double]"0.0008571597498645556"
$ws.Range("P4").Value = [double]"1.352237078121732"
$ws.Range("Q4").Value = "[0.723289599925578, 1.9811845563178858]"
$ws.Range("R4").Value = [double]"8.220665341363009e-05"
$ws.Range("S4").Value = [double]"8.220665341363009e-05"
$ws.Range("T4").Value = [double]"48.82031037082281"
$ws.Range("U4").Value = "[37.49611649569451, 60.14450424595111]"
$ws.Range("V4").Value = [double]"3.570099771366131e-11"
$ws.Range("W4").Value = [double]"3.570099771366131e-11"
$ws.Range("X4").Value = [double]"18.6778778778781"
$ws.Range("Y4").Value = [double]"16.2954954954957"
$ws.Range("Z4").Value = [double]"21.06026026026051"
$ws.Range("F5").Value = [double]"23.80000000000028"
$ws.Range("H5").Value = [double]"0.0005092830219424194"
$ws.Range("I5").Value = [double]"0.0005092830219424194"
$ws.Range("L5").Value = [double]"45.4691094664972"
$ws.Range("M5").Value = "[18.616326564763625, 72.32189236823078]"
$ws.Range("N5").Value = [double]"0.001379391356885806"
$ws.Range("O5").Value = [double]"0.001379391356885806"
$ws.Range("P5").Value = [double]"1.83023716155081"
$ws.Range("Q5").Value = "[1.1258159859711148, 2.5346583371305043]"
$ws.Range("R5").Value = [double]"4.217794747685133e-06"
$ws.Range("S5").Value = [double]"4.217794747685133e-06"
$ws.Range("T5").Value = [double]"52.60207013019895"
$ws.Range("U5").Value = "[37.417519433432716, 67.78662082696518]"
$ws.Range("V5").Value = [double]"1.10140412346027e-08"
$ws.Range("W5").Value = [double]"1.10140412346027e-08"
$ws.Range("X5").Value = [double]"16.86726726726747"
$ws.Range("Y5").Value = [double]"14.19899899899917"
$ws.Range("Z5").Value = [double]"19.53553553553578"
$ws.Range("F6").Value = [double]"23.80000000000028"
$ws.Range("H6").Value = [double]"0.0001812669036813652"
$ws.Range("I6").Value = [double]"0.0001812669036813652"
$ws.Range("L6").Value = [double]"43.98462586017584"
$ws.Range("M6").Value = "[20.637199477788783, 67.3320522425629]"
$ws.Range("N6").Value = [double]"0.0004393994579445781"
$ws.Range("O6").Value = [double]"0.0004393994579445781"
$ws.Range("P6").Value = [double]"1.704447665911579"
$ws.Range("Q6").Value = "[1.0755001877154244, 2.333395144107734]"
$ws.Range("R6").Value = [double]"1.973772644436877e-06"
$ws.Range("S6").Value = [double]"1.973772644436877e-06"
$ws.Range("T6").Value = [double]"66.55755354790627"
$ws.Range("U6").Value = "[52.94379619442904, 80.17131090138349]"
$ws.Range("V6").Value = [double]"8.411049634560186e-13"
$ws.Range("W6").Value = [double]"8.411049634560186e-13"
$ws.Range("X6").Value = [double]"17.34374374374395"
$ws.Range("Y6").Value = [double]"14.96136136136154"
$ws.Range("Z6").Value = [double]"19.72612612612636"
$ws.Range("F7").Value = [double]"23.80000000000028"
$ws.Range("H7").Value = [double]"1.597541977749195e-05"
$ws.Range("I7").Value = [double]"1.597541977749195e-05"
$ws.Range("L7").Value = [double]"41.75857043119863"
$ws.Range("M7").Value = "[20.13639834531414, 63.38074251708313]"
$ws.Range("N7").Value = [double]"0.0003280702460848328"
$ws.Range("O7").Value = [double]"0.0003280702460848328"
$ws.Range("P7").Value = [double]"2.257921446724195"
$ws.Range("Q7").Value = "[1.7547634641672705, 2.7610794292811196]"
$ws.Range("R7").Value = [double]"1.118882764217233e-11"
$ws.Range("S7").Value = [double]"1.118882764217233e-11"
$ws.Range("T7").Value = [double]"60.43840013920259"
$ws.Range("U7").Value = "[49.32531081951519, 71.55148945888999]"
$ws.Range("V7").Value = [double]"2.775557561562891e-14"
$ws.Range("W7").Value = [double]"2.775557561562891e-14"
$ws.Range("X7").Value = [double]"15.24724724724743"
$ws.Range("Y7").Value = [double]"13.3413413413415"
$ws.Range("Z7").Value = [double]"17.15315315315336"
$ws.Range("F8").Value = [double]"23.80000000000028"
$ws.Range("H8").Value = [double]"2.995735828803525e-05"
$ws.Range("I8").Value = [double]"2.995735828803525e-05"
$ws.Range("L8").Value = [double]"43.88200145460301"
$ws.Range("M8").Value = "[24.164637205248887, 63.599365703957126]"
$ws.Range("N8").Value = [double]"5.037655934669871e-05"
$ws.Range("O8").Value = [double]"5.037655934669871e-05"
$ws.Range("Q8").Value = "[1.1132370364071917, 2.220184598032427]"
$ws.Range("R8").Value = [double]"2.498463518296745e-07"
$ws.Range("S8").Value = [double]"2.498463518296745e-07"
$ws.Range("T8").Value = [double]"67.01921134445361"
$ws.Range("U8").Value = "[54.93943735647021, 79.09898533243701]"
$ws.Range("V8").Value = [double]"1.4210854715202e-14"
$ws.Range("W8").Value = [double]"1.4210854715202e-14"
$ws.Range("X8").Value = [double]"17.48668668668689"
$ws.Range("Y8").Value = [double]"15.39019019019037"
$ws.Range("Z8").Value = [double]"19.58318318318342"
$ws.Range("F9").Value = [double]"23.80000000000028"
$ws.Range("H9").Value = [double]"9.743988533261749e-06"
$ws.Range("I9").Value = [double]"9.743988533261749e-06"
$ws.Range("L9").Value = [double]"49.47070792717385"
$ws.Range("M9").Value = "[24.446201552900035, 74.49521430144766]"
$ws.Range("N9").Value = [double]"0.0002469548513202025"
$ws.Range("O9").Value = [double]"0.0002469548513202025"
$ws.Range("P9").Value = [double]"2.207605648468503"
$ws.Range("Q9").Value = "[1.7044476659115793, 2.7107636310254266]"
$ws.Range("R9").Value = [double]"2.156719247636829e-11"
$ws.Range("S9").Value = [double]"2.156719247636829e-11"
$ws.Range("T9").Value = [double]"66.19177790968139"
$ws.Range("U9").Value = "[53.38090259571982, 79.00265322364297]"
$ws.Range("V9").Value = [double]"1.469935284603707e-13"
$ws.Range("W9").Value = [double]"1.469935284603707e-13"
$ws.Range("X9").Value = [double]"15.43783783783802"
$ws.Range("Y9").Value = [double]"13.53193193193209"
$ws.Range("Z9").Value = [double]"17.34374374374395"
$ws.Range("B10").Value = [double]"1"
$ws.Range("F10").Value = [double]"23.80000000000028"
$ws.Range("H10").Value = [double]"0.0004474018792164358"
$ws.Range("I10").Value = [double]"0.0004474018792164358"
$ws.Range("L10").Value = [double]"36.45897232285434"
$ws.Range("M10").Value = "[15.286051160609766, 57.63189348509891]"
$ws.Range("N10").Value = [double]"0.001165323461324475"
$ws.Range("O10").Value = [double]"0.001165323461324475"
$ws.Range("P10").Value = [double]"1.943447707626118"
$ws.Range("Q10").Value = "[1.301921279866039, 2.5849741353861964]"
$ws.Range("R10").Value = [double]"2.206328975074001e-07"
$ws.Range("S10").Value = [double]"2.206328975074001e-07"
$ws.Range("T10").Value = [double]"56.6603233941688"
$ws.Range("U10").Value = "[44.59411493673073, 68.72653185160686]"
$ws.Range("V10").Value = [double]"2.894351425197783e-12"
$ws.Range("W10").Value = [double]"2.894351425197783e-12"
$ws.Range("X10").Value = [double]"16.43843843843863"
$ws.Range("Y10").Value = [double]"14.00840840840857"
$ws.Range("Z10").Value = [double]"18.86846846846869"
$ws.Range("B11").Value = [double]"0"
$ws.Range("F11").Value = [double]"23.16000000000018"
$ws.Range("H11").Value = [double]"0.06166352672390341"
$ws.Range("I11").Value = [double]"0.06166352672390341"
$ws.Range("L11").Value = [double]"19.86971052569707"
$ws.Range("M11").Value = "[-0.8719102250083139, 40.61133127640246]"
$ws.Range("N11").Value = [double]"0.05999637672113933"
$ws.Range("O11").Value = [double]"0.05999637672113933"
$ws.Range("P11").Value = [double]"1.389973926813503"
$ws.Range("Q11").Value = "[-0.018868424345883206, 2.7988162779728887]"
$ws.Range("R11").Value = [double]"0.05301574223404937"
$ws.Range("S11").Value = [double]"0.05301574223404937"
$ws.Range("T11").Value = [double]"38.76209776906323"
$ws.Range("U11").Value = "[27.337519196644557, 50.186676341481906]"
$ws.Range("V11").Value = [double]"1.799441529115597e-08"
$ws.Range("W11").Value = [double]"1.799441529115597e-08"
$ws.Range("X11").Value = [double]"18.03651651651666"
$ws.Range("Y11").Value = [double]"12.84348348348358"
$ws.Range("Z11").Value = [double]"23.22954954954973"
$ws.Range("F12").Value = [double]"23.16000000000018"
$ws.Range("H12").Value = [double]"4.256401239843832e-05"
$ws.Range("I12").Value = [double]"4.256401239843832e-05"
$ws.Range("L12").Value = [double]"42.21886769180235"
$ws.Range("M12").Value = "[23.362505477286774, 61.07522990631793]"
$ws.Range("N12").Value = [double]"4.615415049191718e-05"
$ws.Range("O12").Value = [double]"4.615415049191718e-05"
$ws.Range("P12").Value = [double]"1.540921321580579"
$ws.Range("Q12").Value = "[0.9874475407679633, 2.094395102393195]"
$ws.Range("R12").Value = [double]"1.190069760736989e-06"
$ws.Range("S12").Value = [double]"1.190069760736989e-06"
$ws.Range("T12").Value = [double]"66.18614380567146"
$ws.Range("U12").Value = "[54.45034520988506, 77.92194240145787]"
$ws.Range("V12").Value = [double]"8.215650382226158e-15"
$ws.Range("W12").Value = [double]"8.215650382226158e-15"
$ws.Range("X12").Value = [double]"17.48012012012026"
$ws.Range("Y12").Value = [double]"15.44000000000012"
$ws.Range("Z12").Value = [double]"19.52024024024039"
$ws.Range("F13").Value = [double]"23.16000000000018"
$ws.Range("H13").Value = [double]"4.950507670187232e-05"
$ws.Range("I13").Value = [double]"4.950507670187232e-05"
$ws.Range("L13").Value = [double]"54.07706661653019"
$ws.Range("M13").Value = "[27.440965772018174, 80.7131674610422]"
$ws.Range("N13").Value = [double]"0.0001766038367954081"
$ws.Range("O13").Value = [double]"0.0001766038367954081"
$ws.Range("P13").Value = [double]"1.880552959806503"
$ws.Range("Q13").Value = "[1.2893423303021185, 2.4717635893108882]"
$ws.Range("R13").Value = [double]"7.764234122298319e-08"
$ws.Range("S13").Value = [double]"7.764234122298319e-08"
$ws.Range("T13").Value = [double]"68.10397386922585"
$ws.Range("U13").Value = "[52.87685501816634, 83.33109272028537]"
$ws.Range("V13").Value = [double]"1.233901869568399e-11"
$ws.Range("W13").Value = [double]"1.233901869568399e-11"
$ws.Range("X13").Value = [double]"16.22822822822835"
$ws.Range("Y13").Value = [double]"14.04900900900912"
$ws.Range("Z13").Value = [double]"18.40744744744759"
$ws.Range("F14").Value = [double]"23.16000000000018"
$ws.Range("H14").Value = [double]"1.080018715715614e-05"
$ws.Range("I14").Value = [double]"1.080018715715614e-05"
$ws.Range("L14").Value = [double]"47.06194846160798"
$ws.Range("M14").Value = "[23.914320265208687, 70.20957665800726]"
$ws.Range("N14").Value = [double]"0.0001733911369581342"
$ws.Range("O14").Value = [double]"0.0001733911369581342"
$ws.Range("P14").Value = [double]"2.044079304137503"
$ws.Range("Q14").Value = "[1.515763422452732, 2.5723951858222733]"
$ws.Range("R14").Value = [double]"6.917479922208258e-10"
$ws.Range("S14").Value = [double]"6.917479922208258e-10"
$ws.Range("T14").Value = [double]"70.63352501407005"
$ws.Range("U14").Value = "[58.472079647551226, 82.79497038058886]"
$ws.Range("V14").Value = [double]"3.108624468950438e-15"
$ws.Range("W14").Value = [double]"3.108624468950438e-15"
$ws.Range("X14").Value = [double]"15.62546546546559"
$ws.Range("Y14").Value = [double]"13.67807807807818"
$ws.Range("Z14").Value = [double]"17.572852852853"
